$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 131064784
$ws.Range("B2").Value = 91829
$ws.Range("E2").Value = 5432
$ws.Range("F2").Value = "Granticka"
$ws.Range("G2").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H2").Value = ""
$ws.Range("Q2").Value = 442100
$ws.Range("R2").Value = 7039221
$ws.Range("AC2").Value = ""

# Row 3
$ws.Range("A3").Value = 131064775
$ws.Range("B3").Value = 57884
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("Q3").Value = 442085
$ws.Range("R3").Value = 7039138
$ws.Range("AC3").Value = "Ringhack"

# Row 4
$ws.Range("A4").Value = 131064773
$ws.Range("B4").Value = 57884
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("Q4").Value = 442108
$ws.Range("R4").Value = 7039138
$ws.Range("AC4").Value = "Ringhack äldre"

# Row 9
$ws.Range("A9").Value = 131064770
$ws.Range("B9").Value = 57884
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = "Tretåig hackspett"
$ws.Range("G9").Value = "Picoides tridactylus"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("Q9").Value = 442198
$ws.Range("R9").Value = 7039206
$ws.Range("AC9").Value = "Ringhack färska och äldre"

# Row 10
$ws.Range("A10").Value = 131064778
$ws.Range("B10").Value = 57884
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("Q10").Value = 442145
$ws.Range("R10").Value = 7039101
$ws.Range("AC10").Value = "Ringhack"

# Row 11
$ws.Range("A11").Value = 131064783
$ws.Range("B11").Value = 91829
$ws.Range("E11").Value = 5432
$ws.Range("F11").Value = "Granticka"
$ws.Range("G11").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H11").Value = ""
$ws.Range("Q11").Value = 442292
$ws.Range("R11").Value = 7039182
$ws.Range("AC11").Value = ""

# Row 12
$ws.Range("A12").Value = 131064766
$ws.Range("B12").Value = 57884
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = "Tretåig hackspett"
$ws.Range("G12").Value = "Picoides tridactylus"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("Q12").Value = 442271
$ws.Range("R12").Value = 7039174
$ws.Range("AC12").Value = "Ringhack äldre"

# Row 17
$ws.Range("A17").Value = 131064772
$ws.Range("B17").Value = 57884
$ws.Range("E17").Value = 100109
$ws.Range("F17").Value = "Tretåig hackspett"
$ws.Range("G17").Value = "Picoides tridactylus"
$ws.Range("H17").Value = "(Linnaeus, 1758)"
$ws.Range("Q17").Value = 442099
$ws.Range("R17").Value = 7039220
$ws.Range("AC17").Value = "Bohål ca 3m upp i grantickerötad granhögstubbe Även ett påbörjat på 2m"

# Row 18
$ws.Range("A18").Value = 131064781
$ws.Range("B18").Value = 91805
$ws.Range("E18").Value = 1108
$ws.Range("F18").Value = "Harticka"
$ws.Range("G18").Value = "Pelloporus leporinus"
$ws.Range("H18").Value = "(Fr.) Krieglst."
$ws.Range("Q18").Value = 442200
$ws.Range("R18").Value = 7039150
$ws.Range("AC18").Value = ""

# Row 19
$ws.Range("A19").Value = 131064780
$ws.Range("B19").Value = 91805
$ws.Range("E19").Value = 1108
$ws.Range("F19").Value = "Harticka"
$ws.Range("G19").Value = "Pelloporus leporinus"
$ws.Range("H19").Value = "(Fr.) Krieglst."
$ws.Range("Q19").Value = 442259
$ws.Range("R19").Value = 7039181
$ws.Range("AC19").Value = ""

